$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text in A2 and B2 per commit: changed synthesize_text() descriptions
$ws.Range("A2").Value = "I can write a for loop"
$ws.Range("B2").Value = "That learner has knowledge of the for loop. You should ask the learner complex questions"

# Update active cell selection to B6
$ws.Range("B6").Select()

# Extend the explicit custom width to column 11 (K) to match the default width
$ws.Columns.Item(11).ColumnWidth = 7.71
